$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")

# Relabel the "Currency" row: lowercase the label and trim the trailing
# space off the value, giving the cell a new (green-filled, Normal-based)
# style.
$ws1.Range("A6").Value = "currency"
$ws1.Range("B6").Value = "US Dollar"
$ws1.Range("B6").Style = "Normal"
$ws1.Range("B6").Interior.Color = 5296274

# Move the active tab/selection from the output sheet to the input sheet,
# landing on the row we just edited.
$ws1.Activate()
$ws1.Range("A6:B6").Select()
